# "BPMNs and Update Hours"
# Log 1 hour of work against the "Develop BPMN Diagrams" task (row 8) and
# against the corresponding task in the second table (row 21), matching the
# look/format of the neighboring "hours worked" input cells on each row.
# All the downstream SUM()/running-total formulas and the burndown chart's
# cached series values recompute automatically from this single input edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BURNDOWN CHART")

$ws.Range("K8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value = 1

$ws.Range("O21").Copy()
$ws.Range("P21").PasteSpecial(-4122)
$ws.Range("P21").Value = 1

# Update the sheet's current selection, as recorded the last time it was saved.
$ws.Range("H28").Select()
